$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal TEXT (matches inlineStr in source),
# avoiding Excel's automatic numeric coercion for number-looking strings,
# and avoiding any lingering style/number-format residue on the cell.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue 2 4 '68.846.70'   # D2
Set-TextValue 2 5 '  -0.18%  '   # E2
Set-TextValue 3 4 '3.854.83'   # D3
Set-TextValue 3 5 '  +2.76%  '   # E3
Set-TextValue 4 5 '  -0.08%  '   # E4
Set-TextValue 5 4 '602.77'   # D5
Set-TextValue 5 5 '  +0.10%  '   # E5
Set-TextValue 6 4 '163.09'   # D6
Set-TextValue 6 5 '  -2.57%  '   # E6
Set-TextValue 7 4 '3.851.20'   # D7
Set-TextValue 7 5 '  +2.71%  '   # E7
Set-TextValue 8 5 '  +0.01%  '   # E8
Set-TextValue 9 5 '  -1.52%  '   # E9
Set-TextValue 10 5 '  -0.93%  '   # E10
Set-TextValue 11 4 '6.31'   # D11
Set-TextValue 11 5 '  -2.53%  '   # E11
Set-TextValue 12 5 '  -0.24%  '   # E12
Set-TextValue 13 4 '36.92'   # D13
Set-TextValue 13 5 '  -2.48%  '   # E13
Set-TextValue 14 4 '0.0000244'   # D14
Set-TextValue 14 5 '  -2.01%  '   # E14
Set-TextValue 15 4 '4.487.38'   # D15
Set-TextValue 15 5 '  +2.42%  '   # E15
Set-TextValue 16 4 '3.834.68'   # D16
Set-TextValue 16 5 '  +2.19%  '   # E16
Set-TextValue 17 4 '69.017.47'   # D17
Set-TextValue 17 5 '  -0.03%  '   # E17
Set-TextValue 18 5 '  +2.56%  '   # E18
Set-TextValue 19 5 '  -0.19%  '   # E19
Set-TextValue 20 4 '11.42'   # D20
Set-TextValue 20 5 '  +5.45%  '   # E20
Set-TextValue 21 4 '17.17'   # D21
Set-TextValue 21 5 '  -0.42%  '   # E21
Set-TextValue 22 4 '485.22'   # D22
Set-TextValue 22 5 '  -1.53%  '   # E22
Set-TextValue 23 4 '0.720'   # D23
Set-TextValue 23 5 '  -0.75%  '   # E23
Set-TextValue 24 4 '0.0000158'   # D24
Set-TextValue 24 5 '  +3.03%  '   # E24
Set-TextValue 25 4 '84.08'   # D25
Set-TextValue 25 5 '  -0.83%  '   # E25
Set-TextValue 26 4 '2.26'   # D26
Set-TextValue 26 5 '  -2.09%  '   # E26
Set-TextValue 27 4 '12.10'   # D27
Set-TextValue 27 5 '  -1.92%  '   # E27
Set-TextValue 28 5 '  -1.20%  '   # E28
Set-TextValue 29 5 '  -0.07%  '   # E29
Set-TextValue 30 5 '  -0.85%  '   # E30
Set-TextValue 31 4 '7.94'   # D31
Set-TextValue 31 5 '  -0.80%  '   # E31
Set-TextValue 32 4 '4.004.82'   # D32
Set-TextValue 32 5 '  +2.74%  '   # E32
Set-TextValue 33 4 '2.39'   # D33
Set-TextValue 33 5 '  -3.68%  '   # E33
Set-TextValue 34 4 '32.20'   # D34
Set-TextValue 34 5 '  +2.04%  '   # E34
Set-TextValue 35 4 '3.801.06'   # D35
Set-TextValue 35 5 '  +3.12%  '   # E35
Set-TextValue 36 5 '  -1.72%  '   # E36
Set-TextValue 37 5 '  +1.33%  '   # E37
Set-TextValue 38 5 '  +4.52%  '   # E38
Set-TextValue 39 4 '5.88'   # D39
Set-TextValue 39 5 '  -0.24%  '   # E39
Set-TextValue 40 4 '0.998'   # D40
Set-TextValue 40 5 '  -0.18%  '   # E40
Set-TextValue 41 2 'Bittensor'   # B41
Set-TextValue 41 3 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'   # C41
Set-TextValue 41 4 '443.20'   # D41
Set-TextValue 41 5 '  +3.30%  '   # E41
Set-TextValue 42 2 'TheGraph'   # B42
Set-TextValue 42 3 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'   # C42
Set-TextValue 42 4 '0.318'   # D42
Set-TextValue 42 5 '  -1.89%  '   # E42
Set-TextValue 43 5 '  +0.99%  '   # E43
Set-TextValue 44 4 '48.54'   # D44
Set-TextValue 44 5 '  -0.15%  '   # E44
Set-TextValue 45 5 '  -1.31%  '   # E45
Set-TextValue 46 5 '  +0.02%  '   # E46
Set-TextValue 47 4 '8.40'   # D47
Set-TextValue 47 5 '  -1.00%  '   # E47
Set-TextValue 48 4 '27.20'   # D48
Set-TextValue 48 5 '  +16.25%  '   # E48
Set-TextValue 49 4 '2.832.94'   # D49
Set-TextValue 49 5 '  +1.56%  '   # E49
Set-TextValue 50 4 '142.42'   # D50
Set-TextValue 50 5 '  +0.61%  '   # E50
Set-TextValue 51 4 '0.0356'   # D51
Set-TextValue 51 5 '  +1.00%  '   # E51
